$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The header row shrinks from A1:N1 (14 columns) to A1:L1 (12 columns):
# 학생번호/비밀번호/핸드폰번호/보호자 이름/가족 관계/상세주소 are dropped or
# merged into renamed columns, and the remaining headers are reordered and
# relabeled to the new layout.
$ws.Range("A1").Value = "주민번호"
$ws.Range("B1").Value = "학과코드"
$ws.Range("C1").Value = "이름"
$ws.Range("D1").Value = "영문이름"
$ws.Range("E1").Value = "이메일"
$ws.Range("F1").Value = "우편번호"
$ws.Range("G1").Value = "주소"
$ws.Range("H1").Value = "상세 주소"
$ws.Range("I1").Value = "핸드폰 번호"
$ws.Range("J1").Value = "집 전화번호"
$ws.Range("K1").Value = "보호자이름"
$ws.Range("L1").Value = "가족관계"

# Columns M and N (옛 상세주소/이메일) no longer exist in the header row, so
# the used range shrinks back down to A1:L1.
$ws.Range("M1:N1").Clear()

# Selection moves to M1 (just past the new last header column), matching
# the saved view state.
$ws.Range("M1").Select() | Out-Null
